$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 192.75
$ws.Range("I2").Value = 185.5
$ws.Range("K2").Value = 185.5
$ws.Range("M2").Value = -72.5
$ws.Range("H61").Value = 87.5
$ws.Range("I61").Value = 87.5
$ws.Range("K61").Value = 262.5
$ws.Range("M61").Value = -90.5
$ws.Range("H64").Value = 3462.6
$ws.Range("J64").Value = 3450
$ws.Range("L64").Value = 3450
$ws.Range("N64").Value = -3946
$ws.Range("H67").Value = 3462.6
$ws.Range("J67").Value = 3450
$ws.Range("L67").Value = 3450
$ws.Range("N67").Value = -5166
$ws.Range("H69").Value = 4496.875
$ws.Range("I69").Value = 6750
$ws.Range("J69").Value = 3745.8333
$ws.Range("K69").Value = 20250
$ws.Range("L69").Value = 11237.4999
$ws.Range("M69").Value = -19376
$ws.Range("N69").Value = -12985.4999
$ws.Range("H70").Value = 1370
$ws.Range("I70").Value = 1100
$ws.Range("J70").Value = 1443.6364
$ws.Range("K70").Value = 3300
$ws.Range("L70").Value = 4330.9092
$ws.Range("M70").Value = -3030
$ws.Range("N70").Value = -4870.9092
$ws.Range("H72").Value = 4496.875
$ws.Range("I72").Value = 6750
$ws.Range("J72").Value = 3745.8333
$ws.Range("K72").Value = 60750
$ws.Range("L72").Value = 33712.4997
$ws.Range("M72").Value = -56382
$ws.Range("N72").Value = -42448.4997
$ws.Range("H73").Value = 1370
$ws.Range("I73").Value = 1100
$ws.Range("J73").Value = 1443.6364
$ws.Range("K73").Value = 3300
$ws.Range("L73").Value = 4330.9092
$ws.Range("M73").Value = -2364
$ws.Range("N73").Value = -6202.9092
$ws.Range("H74").Value = 4149.2
$ws.Range("I74").Value = 4246
$ws.Range("J74").Value = 4125
$ws.Range("K74").Value = 4246
$ws.Range("L74").Value = 4125
$ws.Range("M74").Value = -3310
$ws.Range("N74").Value = -5997
$ws.Range("H75").Value = 34663
$ws.Range("J75").Value = 34663
$ws.Range("L75").Value = 34663
$ws.Range("N75").Value = -36535
$ws.Range("H77").Value = 4149.2
$ws.Range("I77").Value = 4246
$ws.Range("J77").Value = 4125
$ws.Range("K77").Value = 21230
$ws.Range("L77").Value = 20625
$ws.Range("M77").Value = -16550
$ws.Range("N77").Value = -29985
$ws.Range("H78").Value = 34663
$ws.Range("J78").Value = 34663
$ws.Range("L78").Value = 103989
$ws.Range("N78").Value = -113349
$ws.Range("H80").Value = 1601.625
$ws.Range("I80").Value = 2566.4
$ws.Range("J80").Value = 1163.091
$ws.Range("K80").Value = 7699.200000000001
$ws.Range("L80").Value = 3489.273
$ws.Range("M80").Value = -6701.200000000001
$ws.Range("N80").Value = -5485.272999999999
$ws.Range("H81").Value = 32000
$ws.Range("J81").Value = 32000
$ws.Range("L81").Value = 32000
$ws.Range("N81").Value = -33996
$ws.Range("H82").Value = 6210.5454
$ws.Range("I82").Value = 2406.3333
$ws.Range("J82").Value = 7637.125
$ws.Range("K82").Value = 7218.999899999999
$ws.Range("L82").Value = 22911.375
$ws.Range("M82").Value = -6812.999899999999
$ws.Range("N82").Value = -23723.375
$ws.Range("H83").Value = 1601.625
$ws.Range("I83").Value = 2566.4
$ws.Range("J83").Value = 1163.091
$ws.Range("K83").Value = 23097.6
$ws.Range("L83").Value = 10467.819
$ws.Range("M83").Value = -18105.6
$ws.Range("N83").Value = -20451.819
$ws.Range("H84").Value = 32000
$ws.Range("J84").Value = 32000
$ws.Range("L84").Value = 96000
$ws.Range("N84").Value = -105984
$ws.Range("H85").Value = 6210.5454
$ws.Range("I85").Value = 2406.3333
$ws.Range("J85").Value = 7637.125
$ws.Range("K85").Value = 7218.999899999999
$ws.Range("L85").Value = 22911.375
$ws.Range("M85").Value = -5814.999899999999
$ws.Range("N85").Value = -25719.375
$ws.Range("H86").Value = 3945.7273
$ws.Range("I86").Value = 3150.5
$ws.Range("K86").Value = 3150.5
$ws.Range("M86").Value = -2027.5
$ws.Range("H88").Value = 10000
$ws.Range("I88").Value = 10000
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 10000
$ws.Range("L88").Value = 0
$ws.Range("M88").Value = -9594
$ws.Range("N88").ClearContents()
$ws.Range("H89").Value = 3945.7273
$ws.Range("I89").Value = 3150.5
$ws.Range("K89").Value = 15752.5
$ws.Range("M89").Value = -10136.5
$ws.Range("H91").Value = 10000
$ws.Range("I91").Value = 10000
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 10000
$ws.Range("L91").Value = 0
$ws.Range("M91").Value = -8596
$ws.Range("N91").ClearContents()
$ws.Range("H138").Value = 2299.611
$ws.Range("I138").Value = 1758.7222
$ws.Range("J138").Value = 2840.5
$ws.Range("K138").Value = 5276.1666
$ws.Range("L138").Value = 8521.5
$ws.Range("M138").Value = -136.1665999999996
$ws.Range("N138").Value = -18801.5
$ws.Range("H139").Value = 41419.5
$ws.Range("J139").Value = 41419.5
$ws.Range("L139").Value = 41419.5
$ws.Range("N139").Value = -51699.5
$ws.Range("H140").Value = 52925
$ws.Range("J140").Value = 52925
$ws.Range("L140").Value = 52925
$ws.Range("N140").Value = -63285

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H59").Value = 29019.666
$ws.Range("J59").Value = 29019.666
$ws.Range("L59").Value = 29019.666
$ws.Range("N59").Value = -30627.666
$ws.Range("H88").Value = 2625.2144
$ws.Range("I88").Value = 2794.3333
$ws.Range("J88").Value = 2498.375
$ws.Range("K88").Value = 2794.3333
$ws.Range("L88").Value = 2498.375
$ws.Range("M88").Value = -2388.3333
$ws.Range("N88").Value = -3310.375
$ws.Range("H91").Value = 2625.2144
$ws.Range("I91").Value = 2794.3333
$ws.Range("J91").Value = 2498.375
$ws.Range("K91").Value = 2794.3333
$ws.Range("L91").Value = 2498.375
$ws.Range("M91").Value = -1390.3333
$ws.Range("N91").Value = -5306.375

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1818.4482
$ws.Range("I86").Value = 1755.1578
$ws.Range("J86").Value = 1938.7
$ws.Range("K86").Value = 1755.1578
$ws.Range("L86").Value = 1938.7
$ws.Range("M86").Value = -632.1578
$ws.Range("N86").Value = -4184.7
$ws.Range("H89").Value = 1818.4482
$ws.Range("I89").Value = 1755.1578
$ws.Range("J89").Value = 1938.7
$ws.Range("K89").Value = 8775.789000000001
$ws.Range("L89").Value = 9693.5
$ws.Range("M89").Value = -3159.789000000001
$ws.Range("N89").Value = -20925.5
$ws.Range("H105").Value = 4269.3887
$ws.Range("I105").Value = 3181.7273
$ws.Range("J105").Value = 4747.96
$ws.Range("K105").Value = 3181.7273
$ws.Range("L105").Value = 4747.96
$ws.Range("M105").Value = -1434.7273
$ws.Range("N105").Value = -8241.959999999999

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2373.4119
$ws.Range("I134").Value = 2370.6667
$ws.Range("J134").Value = 2380
$ws.Range("K134").Value = 7112.000100000001
$ws.Range("L134").Value = 7140
$ws.Range("M134").Value = -4577.000100000001
$ws.Range("N134").Value = -12210
$ws.Range("H140").Value = 30220
$ws.Range("J140").Value = 30220
$ws.Range("L140").Value = 30220
$ws.Range("N140").Value = -40580

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H82").Value = 5571.364
$ws.Range("J82").Value = 6039.5
$ws.Range("L82").Value = 18118.5
$ws.Range("N82").Value = -18930.5
$ws.Range("H85").Value = 5571.364
$ws.Range("J85").Value = 6039.5
$ws.Range("L85").Value = 18118.5
$ws.Range("N85").Value = -20926.5
$ws.Range("H113").Value = 1600.8823
$ws.Range("I113").Value = 472.6
$ws.Range("J113").Value = 3212.7144
$ws.Range("K113").Value = 1417.8
$ws.Range("L113").Value = 9638.143199999999
$ws.Range("M113").Value = 752.1999999999998
$ws.Range("N113").Value = -13978.1432
$ws.Range("H130").Value = 9060
$ws.Range("I130").Value = 830
$ws.Range("J130").Value = 10888.889
$ws.Range("K130").Value = 2490
$ws.Range("L130").Value = 32666.667
$ws.Range("M130").Value = 2530
$ws.Range("N130").Value = -42706.667
$ws.Range("H131").Value = 816.4343
$ws.Range("I131").Value = 451.66666
$ws.Range("J131").Value = 839.9677
$ws.Range("K131").Value = 1354.99998
$ws.Range("L131").Value = 2519.9031
$ws.Range("M131").Value = 3685.00002
$ws.Range("N131").Value = -12599.9031

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("M53").ClearContents()
$ws.Range("N53").ClearContents()
$ws.Range("H58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").ClearContents()
$ws.Range("H80").Value = 2994.389
$ws.Range("I80").Value = 2966.6667
$ws.Range("J80").Value = 3133
$ws.Range("K80").Value = 2966.6667
$ws.Range("L80").Value = 3133
$ws.Range("M80").Value = -1968.6667
$ws.Range("N80").Value = -5129
$ws.Range("H83").Value = 2994.389
$ws.Range("I83").Value = 2966.6667
$ws.Range("J83").Value = 3133
$ws.Range("K83").Value = 14833.3335
$ws.Range("L83").Value = 15665
$ws.Range("M83").Value = -9841.333500000001
$ws.Range("N83").Value = -25649
$ws.Range("H138").Value = 58599.5
$ws.Range("J138").Value = 58599.5
$ws.Range("L138").Value = 58599.5
$ws.Range("N138").Value = -68879.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()
$ws.Range("H139").Value = 46429.4
$ws.Range("J139").Value = 47071.555
$ws.Range("L139").Value = 47071.555
$ws.Range("N139").Value = -57351.555
